$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 33 data)
$ws.Range("D2").Value = 45237
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 22000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 22000
$ws.Range("S2").Value = 2200

# Row 3 (was row 34 data)
$ws.Range("D3").Value = 45237
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("S3").Value = 1800

# Row 4 (was row 30 data)
$ws.Range("D4").Value = 45224
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 2000

# Row 5 (was row 17 data)
$ws.Range("D5").Value = 44447
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 2150

# Row 6 (was row 22 data)
$ws.Range("D6").Value = 45203
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 21000
$ws.Range("S6").Value = 2100

# Row 7 (was row 13 data)
$ws.Range("D7").Value = 45205
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 22000
$ws.Range("S7").Value = 2200

# Row 8 (was row 38 data)
$ws.Range("D8").Value = 44848
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 24500
$ws.Range("S8").Value = 2450

# Row 9 (was row 39 data)
$ws.Range("D9").Value = 44848
$ws.Range("M9").Value = 120

# Row 10 (was row 20 data)
$ws.Range("D10").Value = 44448
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 21000
$ws.Range("P10").Value = 21500
$ws.Range("S10").Value = 2150

# Row 11 (was row 41 data)
$ws.Range("D11").Value = 45180
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 22000
$ws.Range("S11").Value = 2200

# Row 12 (was row 14 data)
$ws.Range("D12").Value = 44460
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 31000
$ws.Range("O12").Value = 32000
$ws.Range("P12").Value = 31500
$ws.Range("S12").Value = 3150

# Row 13 (was row 15 data)
$ws.Range("D13").Value = 44460
$ws.Range("N13").Value = 30000
$ws.Range("O13").Value = 30000
$ws.Range("P13").Value = 30000
$ws.Range("S13").Value = 3000

# Row 14 (was row 23 data)
$ws.Range("D14").Value = 44874
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 25000
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 25000
$ws.Range("S14").Value = 2500

# Row 15 (was row 24 data)
$ws.Range("D15").Value = 44874
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 23000
$ws.Range("O15").Value = 24000
$ws.Range("P15").Value = 23500
$ws.Range("S15").Value = 2350

# Row 16 (was row 40 data)
$ws.Range("D16").Value = 45243
$ws.Range("M16").Value = 60

# Row 17 (was row 16 data)
$ws.Range("D17").Value = 45191
$ws.Range("M17").Value = 30
$ws.Range("O17").Value = 21000
$ws.Range("P17").Value = 21000
$ws.Range("S17").Value = 2100

# Row 18 (was row 28 data)
$ws.Range("D18").Value = 45189
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 22000
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 22000
$ws.Range("S18").Value = 2200

# Row 19 (was row 5 data)
$ws.Range("D19").Value = 45247
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 20000
$ws.Range("S19").Value = 2000

# Row 20 (was row 31 data)
$ws.Range("D20").Value = 44839
$ws.Range("M20").Value = 120
$ws.Range("N20").Value = 25000
$ws.Range("O20").Value = 26000
$ws.Range("P20").Value = 25500
$ws.Range("S20").Value = 2550

# Row 21 (was row 35 data)
$ws.Range("D21").Value = 45225
$ws.Range("M21").Value = 80
$ws.Range("N21").Value = 21000
$ws.Range("O21").Value = 21000
$ws.Range("P21").Value = 21000
$ws.Range("S21").Value = 2100

# Row 22 (was row 8 data)
$ws.Range("D22").Value = 45212
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 22000
$ws.Range("S22").Value = 2200

# Row 23 (was row 4 data)
$ws.Range("D23").Value = 44868
$ws.Range("M23").Value = 60
$ws.Range("N23").Value = 26000
$ws.Range("O23").Value = 26000
$ws.Range("P23").Value = 26000
$ws.Range("S23").Value = 2600

# Row 24 (was row 3 data)
$ws.Range("D24").Value = 45194
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 22000
$ws.Range("P24").Value = 22000
$ws.Range("S24").Value = 2200

# Row 25 (was row 2 data)
$ws.Range("D25").Value = 44487
$ws.Range("M25").Value = 30
$ws.Range("O25").Value = 24000
$ws.Range("P25").Value = 23500
$ws.Range("S25").Value = 2350

# Row 26 (was row 9 data)
$ws.Range("D26").Value = 44446

# Row 27 (was row 7 data)
$ws.Range("D27").Value = 45216
$ws.Range("M27").Value = 60

# Row 28 (was row 6 data)
$ws.Range("D28").Value = 45236
$ws.Range("M28").Value = 100

# Row 29 (was row 25 data)
$ws.Range("D29").Value = 45196
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 23000
$ws.Range("O29").Value = 23000
$ws.Range("P29").Value = 23000
$ws.Range("S29").Value = 2300

# Row 30 (was row 27 data)
$ws.Range("D30").Value = 45230
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 21000
$ws.Range("O30").Value = 21000
$ws.Range("P30").Value = 21000
$ws.Range("S30").Value = 2100

# Row 31 (was row 26 data)
$ws.Range("D31").Value = 44452
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 21000
$ws.Range("O31").Value = 22000
$ws.Range("P31").Value = 21500
$ws.Range("S31").Value = 2150

# Row 32 (was row 36 data)
$ws.Range("D32").Value = 45217
$ws.Range("M32").Value = 30
$ws.Range("N32").Value = 21000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 21000
$ws.Range("S32").Value = 2100

# Row 33 (was row 12 data)
$ws.Range("D33").Value = 45176
$ws.Range("M33").Value = 30

# Row 34 (was row 11 data)
$ws.Range("D34").Value = 44841
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 60
$ws.Range("N34").Value = 23000
$ws.Range("O34").Value = 24000
$ws.Range("P34").Value = 23500
$ws.Range("S34").Value = 2350

# Row 35 (was row 37 data)
$ws.Range("D35").Value = 45244
$ws.Range("M35").Value = 90
$ws.Range("N35").Value = 20000
$ws.Range("P35").Value = 20444
$ws.Range("S35").Value = 2044

# Row 36 (was row 10 data)
$ws.Range("D36").Value = 45209
$ws.Range("M36").Value = 50
$ws.Range("N36").Value = 22000
$ws.Range("O36").Value = 22000
$ws.Range("P36").Value = 22000
$ws.Range("S36").Value = 2200

# Row 37 (was row 29 data)
$ws.Range("D37").Value = 45219
$ws.Range("M37").Value = 30
$ws.Range("O37").Value = 20000
$ws.Range("P37").Value = 20000
$ws.Range("S37").Value = 2000

# Row 38 (was row 18 data)
$ws.Range("D38").Value = 44461
$ws.Range("N38").Value = 31000
$ws.Range("O38").Value = 32000
$ws.Range("P38").Value = 31500
$ws.Range("S38").Value = 3150

# Row 39 (was row 19 data)
$ws.Range("D39").Value = 44461
$ws.Range("M39").Value = 30
$ws.Range("N39").Value = 30000
$ws.Range("O39").Value = 30000
$ws.Range("P39").Value = 30000
$ws.Range("S39").Value = 3000

# Row 40 (was row 32 data)
$ws.Range("D40").Value = 45239
$ws.Range("M40").Value = 100
$ws.Range("N40").Value = 22000
$ws.Range("O40").Value = 22000
$ws.Range("P40").Value = 22000
$ws.Range("S40").Value = 2200

# Row 41 (was row 21 data)
$ws.Range("D41").Value = 45173
$ws.Range("M41").Value = 50

